# Update "三次产业固定资产（不含农户）投资额及增速" sheet:
#  - remove the oldest 7 years of data (2003-2009)
#  - keep 2010-2020 rows (values unchanged, only shift up after deletion)
#  - append two new years: 2021 and 2022

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 8 (years 2003-2009), shifting rows 9-19 up to rows 2-12
$ws.Range("A2:G8").EntireRow.Delete() | Out-Null

# Existing data now occupies rows 2-12 (years 2010-2020); append new rows for 2021 and 2022
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 14275
$ws.Range("C13").Value = 9.1
$ws.Range("D13").Value = 362877
$ws.Range("E13").Value = 2.1
$ws.Range("F13").Value = 167395
$ws.Range("G13").Value = 11.3

$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 14293
$ws.Range("C14").Value = 0.2
$ws.Range("D14").Value = 373842
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 184004
$ws.Range("G14").Value = 10.3

# Match the style used for the A column year labels (centered, bordered cell style)
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13:A14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
